# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 234 (pushing the existing
# rows 234-255 down to 235-256), with the new dimension A1:T256.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 234, shifting rows 234-255 down
# to 235-256 (values, formats and all).
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with this week's record.
$ws.Range("A234").Value = 10
$ws.Range("B234").Value = "Vega Modelo de Temuco"
$ws.Range("C234").Value = "La Araucanía"
$ws.Range("D234").Value = 45106
$ws.Range("E234").Value = 9
$ws.Range("F234").Value = "Fruta"
$ws.Range("G234").Value = 100104
$ws.Range("H234").Value = "Frutos de pepita"
$ws.Range("I234").Value = 100104001
$ws.Range("J234").Value = "Granada"
$ws.Range("K234").Value = "Wonderfull"
$ws.Range("L234").Value = "Primera"
$ws.Range("M234").Value = 180
$ws.Range("N234").Value = 12000
$ws.Range("O234").Value = 13000
$ws.Range("P234").Value = 12556
$ws.Range("Q234").Value = '$/bandeja 10 kilos'
$ws.Range("R234").Value = "Provincia de Limarí"
$ws.Range("S234").Value = 1256
$ws.Range("T234").Value = 10
